# Update "horarios" workbook with the latest scrape (04:19:04) results.
$wb = $excel.ActiveWorkbook

$updateTime = "04:19:04"

# ----------------------------------------------------------------
# Sheet 1: LP1912
# ----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value2 = "Última actualización: $updateTime"
$ws1.Range("A3").Value2 = "Total filas: 16"

$ws1.Cells.Item(19, 1).Value2 = $updateTime
$ws1.Cells.Item(19, 2).Value2 = "06:04"
$ws1.Cells.Item(19, 3).Value2 = "16_SANTA ANA"
$ws1.Cells.Item(19, 4).Value2 = 105
$ws1.Cells.Item(19, 5).Value2 = "LP1912"

$ws1.Cells.Item(20, 1).Value2 = $updateTime
$ws1.Cells.Item(20, 2).Value2 = "06:12"
$ws1.Cells.Item(20, 3).Value2 = "215A_EL PATO"
$ws1.Cells.Item(20, 4).Value2 = 113
$ws1.Cells.Item(20, 5).Value2 = "LP1912"

$ws1.Cells.Item(21, 1).Value2 = $updateTime
$ws1.Cells.Item(21, 2).Value2 = "06:14"
$ws1.Cells.Item(21, 3).Value2 = "225_HARAS DEL SUR"
$ws1.Cells.Item(21, 4).Value2 = 115
$ws1.Cells.Item(21, 5).Value2 = "LP1912"

# ----------------------------------------------------------------
# Sheet 2: LP1912-215
# ----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value2 = "Última actualización: $updateTime"
$ws2.Range("A3").Value2 = "Total filas: 6"

$ws2.Cells.Item(11, 1).Value2 = $updateTime
$ws2.Cells.Item(11, 2).Value2 = "06:12"
$ws2.Cells.Item(11, 3).Value2 = "215A_EL PATO"
$ws2.Cells.Item(11, 4).Value2 = 113
$ws2.Cells.Item(11, 5).Value2 = "LP1912"

# ----------------------------------------------------------------
# Sheet 3: 6203-6173
# ----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value2 = "Última actualización: $updateTime"
$ws3.Range("A3").Value2 = "Total filas: 3"

$ws3.Cells.Item(8, 1).Value2 = $updateTime
$ws3.Cells.Item(8, 2).Value2 = "06:09"
$ws3.Cells.Item(8, 3).Value2 = "215A_LA PLATA"
$ws3.Cells.Item(8, 4).Value2 = 110
$ws3.Cells.Item(8, 5).Value2 = "L6173"
